$d = $word.ActiveDocument

$replacements = @(
    @("98÷7=", "75÷4="),
    @("76÷9=", "39÷8="),
    @("32÷2=", "70÷9="),
    @("30÷3=", "99÷5="),
    @("78÷2=", "70÷5="),
    @("71÷5=", "83÷5="),
    @("40÷3=", "11÷5="),
    @("35÷7=", "13÷4="),
    @("26÷8=", "40÷7="),
    @("83÷9=", "37÷7="),
    @("97÷8=", "73÷7="),
    @("67÷9=", "44÷8="),
    @("99÷9=", "65÷9="),
    @("25÷4=", "39÷5="),
    @("58÷6=", "73÷4="),
    @("41÷8=", "70÷8="),
    @("49÷5=", "58÷2="),
    @("94÷5=", "60÷4="),
    @("53÷2=", "67÷3="),
    @("74÷8=", "76÷3="),
    @("21÷8=", "17÷6="),
    @("51÷6=", "74÷2="),
    @("12÷6=", "74÷2="),
    @("53÷3=", "11÷5="),
    @("87÷8=", "89÷3=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $new, 2)
}
